# Update "想去人数" (want-to-go count) figures that were refreshed by the
# gh-pages data generation run (commit 456a3b4).
#
# Sheet "展览" and sheet "全部类型" both list the same rows of data, and on
# both sheets the following rows need their column F value bumped by 1:
#   F3:  3053 -> 3054
#   F7:  1670 -> 1671
#   F12: 1373 -> 1374
#   F18: 75   -> 76

$wb = $excel.ActiveWorkbook

$targetSheets = @("展览", "全部类型")
$updates = @{
    "F3"  = 3054
    "F7"  = 1671
    "F12" = 1374
    "F18" = 76
}

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
